$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the newly-added simulation-result columns (C/D/E for the first
# few rows, F/G/H for the rest of the data rows 3-37) ---------------------

# Row 3
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0

# Row 4
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0

# Row 5
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 1

# Row 6
$ws.Range("C6").Value = 6
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 2
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0

# Rows 7-37: only F/G/H are newly populated
$ws.Range("F7").Value = 3
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 1

$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0

$ws.Range("F9").Value = 2
$ws.Range("G9").Value = 2
$ws.Range("H9").Value = 1

$ws.Range("F10").Value = 3
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0

$ws.Range("F11").Value = 3
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = 1

$ws.Range("F12").Value = 3
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 1

$ws.Range("F13").Value = 3
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = 1

$ws.Range("F14").Value = 5
$ws.Range("G14").Value = 2
$ws.Range("H14").Value = 2

$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 1

$ws.Range("F16").Value = 2
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 1

$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 1

$ws.Range("F18").Value = 2
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = 0

$ws.Range("F19").Value = 3
$ws.Range("G19").Value = 1
$ws.Range("H19").Value = 0

$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 0

$ws.Range("F21").Value = 2
$ws.Range("G21").Value = 1
$ws.Range("H21").Value = 0

$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 1

$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 1

$ws.Range("F24").Value = 2
$ws.Range("G24").Value = 1
$ws.Range("H24").Value = 1

$ws.Range("F25").Value = 2
$ws.Range("G25").Value = 1
$ws.Range("H25").Value = 0

$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 0

$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 1

$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 0

$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 1

$ws.Range("F30").Value = 5
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = 0

$ws.Range("F31").Value = 3
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 0

$ws.Range("F32").Value = 2
$ws.Range("G32").Value = 1
$ws.Range("H32").Value = 1

$ws.Range("F33").Value = 3
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 0

$ws.Range("F34").Value = 5
$ws.Range("G34").Value = 2
$ws.Range("H34").Value = 1

$ws.Range("F35").Value = 1
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 0

$ws.Range("F36").Value = 4
$ws.Range("G36").Value = 3
$ws.Range("H36").Value = 3

$ws.Range("F37").Value = 4
$ws.Range("G37").Value = 3
$ws.Range("H37").Value = 2

# --- Row 38: drop the old "SUM" label cell and extend the SUM() formulas
# that used to start at column I so they also cover the new F/G/H columns --
$ws.Range("B38").ClearContents()
$ws.Range("F38:K38").ClearContents()
$ws.Range("F38:K38").Formula = "=SUM(F3:F37)"

# --- Update the view: scroll so column G is the leftmost visible column and
# select D7 (matches the new sheetView saved in the workbook) -------------
$ws.Range("D7").Select()
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 1

Write-Output "edit applied"
